$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
# Selection range changed from A7:XFD15 to A7:XFD14 (activeCell stays A7)
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Activate()
$wsSummary.Range("A7:XFD14").Select()

# --- Sheet: Repayment schedule ---
# A new column O (with values matching the existing pattern) is inserted
# between N and P for rows 2-15. Row 2 and 4 (the disbursement rows) are
# left blank like their neighbouring M/N cells; all other rows get 0.
$wsRepay = $wb.Worksheets.Item("Repayment schedule")

$oValues = @{
    2  = $null
    3  = 0
    4  = $null
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
}

foreach ($row in 2..15) {
    $val = $oValues[$row]
    $oCell = $wsRepay.Cells.Item($row, 15)
    if ($null -ne $val) {
        $oCell.Value = $val
    }
    # Copy formatting from the neighbouring N cell (style used throughout
    # the sheet for these columns) onto the new O cell.
    $nCell = $wsRepay.Cells.Item($row, 14)
    $nCell.Copy()
    $oCell.PasteSpecial(-4122)
}

# --- Sheet: Transactions ---
# The transaction IDs have been renumbered and the active selection moved.
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Activate()
$wsTrans.Cells.Item(2, 1).Value = 68
$wsTrans.Cells.Item(3, 1).Value = 67
$wsTrans.Cells.Item(4, 1).Value = 66
$wsTrans.Range("A2:L4").Select()
